# Update the "biochar_land" worksheet: rename the bio-oil yield labels to
# biochar yield labels (the plant now tracks biochar output, not bio-oil),
# and update the production/cost table headers accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("biochar_land")

# Feedstock-specific yield labels (rows 13-17), col B
$ws.Range("B13").Value = "beef biochar yield"
$ws.Range("B14").Value = "dairy biochar yield"
$ws.Range("B15").Value = "goat biochar yield"
$ws.Range("B16").Value = "pork biochar yield"
$ws.Range("B17").Value = "poultry biochar yield"

# Table header row 21: production & unit-cost columns now refer to biochar,
# not the generic "input"
$ws.Range("F21").Value = "Production estimates (kg biochar per year)"
$ws.Range("G21").Value = "Unit cost  ($1975/kg biochar)"

# Re-enter the TO EUAW ratio formula across the whole results table so it is
# stored as one filled-down (shared) formula, same as the other columns
$ws.Range("G22:G27").Formula = "=E22/F22"

# Reflect the author's final selection on the sheet
[void]$ws.Range("B15").Select()
